# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with refreshed quote data.
#
# Note: several Price values look like plain decimal numbers (e.g.
# "215.64"). Assigning such a string straight to .Value would make Excel
# auto-convert the cell to a numeric type (and lose exact text
# formatting/precision), whereas the source workbook stores these as
# text. A leading apostrophe forces Excel to keep the input as text,
# matching the original cell's text representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.962.78"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.642.49"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").Value = "'215.64"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "'0.5086"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "'0.2564"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.06388"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").Value = "'19.52"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "'0.07794"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "'4.299"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "1.652.67"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "'0.5477"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "0.0₅7864"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "'64.39"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "26.021.84"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "'1.005"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "'198.49"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").Value = "'4.452"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "'9.979"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "'6.067"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").Value = "'1.006"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'1.874"
$ws.Range("E24").Value = "  -2.91%  "
$ws.Range("D25").Value = "'141.29"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "'0.1156"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("D27").Value = "'6.885"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("D28").Value = "'15.77"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "'1.242"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'0.05036"
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("D31").Value = "'3.265"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'3.195"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "'1.544"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'2.364"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "'0.9003"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").Value = "'2.587"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").Value = "1.136.61"
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("D38").Value = "'0.5510"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("E39").Value = "  +15.61%  "
$ws.Range("D40").Value = "'0.01563"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "'1.005"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "'2.544"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").Value = "'5.623"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "'0.8185"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").Value = "'100.26"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "1.779.67"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "'0.4535"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "'54.98"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("E51").Value = "  +0.41%  "
